$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.225.80"
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("D3").Value = "2.664.25"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.29%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.545"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "2.663.29"
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("E10").Value = "  +0.97%  "
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("E13").Value = "  -1.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.19%  "
$ws.Range("D15").Value = "3.151.67"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("E16").Value = "  -2.63%  "
$ws.Range("D17").Value = "67.224.53"
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("D18").Value = "2.674.86"
$ws.Range("E18").Value = "  -1.76%  "
$ws.Range("E19").Value = "  -1.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "362.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("E22").Value = "  -3.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.19%  "
$ws.Range("E24").Value = "  -4.71%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.92%  "
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("E29").Value = "  -3.54%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "554.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.19%  "
$ws.Range("E32").Value = "  -3.05%  "
$ws.Range("E33").Value = "  -4.53%  "
$ws.Range("E34").Value = "  -1.64%  "
$ws.Range("E35").Value = "  -2.61%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").Value = "  -5.77%  "
$ws.Range("E38").Value = "  -1.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "156.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.48%  "
$ws.Range("E40").Value = "  -2.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.27%  "
$ws.Range("E42").Value = "  -4.50%  "
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.21%  "
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("E47").Value = "  -6.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.588"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "152.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.23%  "
$ws.Range("E50").Value = "  -3.24%  "
$ws.Range("E51").Value = "  -3.39%  "
